$d = $word.ActiveDocument
$findStart = $d.Content.Duplicate
$null = $findStart.Find.Execute("{#show_photo}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPara = $findStart.Paragraphs(1)

$findEnd = $d.Content.Duplicate
$null = $findEnd.Find.Execute("{/show_border}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPara = $findEnd.Paragraphs(1)

$blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$insertPoint = $blockRange.Start
$null = $blockRange.Delete()

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D639D1" w:rsidRDefault="00C663B4" w:rsidP="00340EB6"><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>{#</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>p</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>hoto}PHOTO</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>: in attachment{/</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>photo}</w:t></w:r></w:p><w:p w:rsidR="00340EB6" w:rsidRPr="00790FC3" w:rsidRDefault="00340EB6" w:rsidP="00340EB6"><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">CRIMINAL RECORD: </w:t></w:r><w:r w:rsidR="005C2EC9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="005C2EC9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>criminal_</w:t></w:r><w:r w:rsidR="00790FC3"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>records</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00790FC3"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p w:rsidR="00340EB6" w:rsidRPr="0025593C" w:rsidRDefault="00C663B4" w:rsidP="00340EB6"><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00C663B4"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>{#</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>border}BORDER</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> RECORD: in attachment with t</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>he legend for translation{/</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>border}</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target = $d.Range($insertPoint, $insertPoint)
$null = $target.InsertXML($xml)

$afterRange = $d.Range($insertPoint, $d.Content.End)
$dummyPara = $afterRange.Paragraphs(4)
$null = $dummyPara.Range.Delete()
